# Added test data for Switzerland market.
#
# The workbook currently has three market sheets (Germany, Belgium, Czech)
# that all share the same layout/template. We add a fourth ("Swiss") by
# duplicating the Czech sheet (so its column widths, styles and merged
# cells come along for free) and then updating only the two market-specific
# cells. Finally the selections on the sheets touched by the edit are
# normalised and the new sheet is left as the active tab/sheet.

$wb = $excel.ActiveWorkbook

# Duplicate the Czech sheet, placing the copy right after it; this becomes
# the new "Swiss" sheet.
$czech = $wb.Worksheets.Item("Czech")
[void]$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Fill in the Switzerland-specific values (same cells that differ between
# the other market sheets).
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2344"

# Put the cursor on A9 for the new sheet.
[void]$swiss.Range("A9").Select()

# Normalise the selections left behind on Germany and Czech.
$germany = $wb.Worksheets.Item("Germany")
[void]$germany.Range("A1:XFD1048576").Select()

[void]$czech.Range("A1:XFD1048576").Select()

# Make Swiss the active sheet/tab.
[void]$swiss.Activate()
